$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.00795755968169761
$ws.Range("C2").Value = 0.957559681697613
$ws.Range("D2").Value = 0.0477453580901857
$ws.Range("E2").Value = 0.986737400530504
$ws.Range("F2").Value = 0.989389920424403
$ws.Range("G2").Value = 0.00795755968169761
$ws.Range("H2").Value = 0.986737400530504
$ws.Range("I2").Value = 0.0318302387267905
$ws.Range("J2").Value = 0.954907161803713
$ws.Range("K2").Value = 0.0026525198938992
$ws.Range("L2").Value = 0.00530503978779841
$ws.Range("M2").Value = 0.0026525198938992
$ws.Range("O2").Value = 0.00530503978779841
$ws.Range("P2").Value = 0.0026525198938992
$ws.Range("U2").Value = 0.013262599469496
$ws.Range("V2").Value = 0.00530503978779841
$ws.Range("W2").Value = 0.0026525198938992
$ws.Range("X2").Value = 0.00530503978779841

# Row 3
$ws.Range("B3").Value = 0.00795755968169761
$ws.Range("D3").Value = 0.952254641909814
$ws.Range("E3").Value = 0.00795755968169761
$ws.Range("H3").Value = 0.0026525198938992
$ws.Range("I3").Value = 0.13262599469496
$ws.Range("K3").Value = 0.0026525198938992
$ws.Range("L3").Value = 0.891246684350133
$ws.Range("M3").Value = 0.994694960212202
$ws.Range("N3").Value = 0.00530503978779841
$ws.Range("O3").Value = 0.00795755968169761
$ws.Range("R3").Value = 0.111405835543767
$ws.Range("S3").Value = 0.809018567639257
$ws.Range("V3").Value = 0.00530503978779841
$ws.Range("W3").Value = 0.0026525198938992
$ws.Range("X3").Value = 0.00795755968169761

# Row 4
$ws.Range("B4").Value = 0.981432360742706
$ws.Range("C4").Value = 0.0318302387267905
$ws.Range("F4").Value = 0.00795755968169761
$ws.Range("G4").Value = 0.989389920424403
$ws.Range("H4").Value = 0.00795755968169761
$ws.Range("I4").Value = 0.185676392572944
$ws.Range("J4").Value = 0.0397877984084881
$ws.Range("K4").Value = 0.00795755968169761
$ws.Range("L4").Value = 0.0026525198938992
$ws.Range("N4").Value = 0.0026525198938992
$ws.Range("P4").Value = 0.997347480106101
$ws.Range("R4").Value = 0.0026525198938992
$ws.Range("U4").Value = 0.978779840848806
$ws.Range("V4").Value = 0.00530503978779841
$ws.Range("W4").Value = 0.994694960212202
$ws.Range("X4").Value = 0.984084880636605

# Row 5
$ws.Range("B5").Value = 0.0026525198938992
$ws.Range("C5").Value = 0.00795755968169761
$ws.Range("E5").Value = 0.00530503978779841
$ws.Range("F5").Value = 0.0026525198938992
$ws.Range("G5").Value = 0.0026525198938992
$ws.Range("H5").Value = 0.0026525198938992
$ws.Range("I5").Value = 0.649867374005305
$ws.Range("J5").Value = 0.00530503978779841
$ws.Range("K5").Value = 0.986737400530504
$ws.Range("L5").Value = 0.10079575596817
$ws.Range("M5").Value = 0.0026525198938992
$ws.Range("N5").Value = 0.992042440318302
$ws.Range("O5").Value = 0.986737400530504
$ws.Range("R5").Value = 0.885941644562334
$ws.Range("S5").Value = 0.190981432360743
$ws.Range("U5").Value = 0.00795755968169761
$ws.Range("V5").Value = 0.984084880636605
$ws.Range("X5").Value = 0.0026525198938992
